$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Content.Find.Execute("2023-10-13 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-10-14 Saturday", 2) | Out-Null

# Update the division problems in the table. The table has 20 rows x 5
# columns, with data only in rows 1, 5, 9, 13, 17 (1-based). Because some
# values (e.g. "85÷5=") are duplicated across cells, address each cell
# directly by (row, column) rather than relying on a global text replace.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "74÷9="
$t.Cell(1, 2).Range.Text = "54÷2="
$t.Cell(1, 3).Range.Text = "92÷7="
$t.Cell(1, 4).Range.Text = "17÷5="
$t.Cell(1, 5).Range.Text = "91÷2="

$t.Cell(5, 1).Range.Text = "95÷5="
$t.Cell(5, 2).Range.Text = "41÷2="
$t.Cell(5, 3).Range.Text = "30÷9="
$t.Cell(5, 4).Range.Text = "62÷9="
$t.Cell(5, 5).Range.Text = "57÷9="

$t.Cell(9, 1).Range.Text = "71÷7="
$t.Cell(9, 2).Range.Text = "26÷8="
$t.Cell(9, 3).Range.Text = "15÷4="
$t.Cell(9, 4).Range.Text = "60÷4="
$t.Cell(9, 5).Range.Text = "35÷4="

$t.Cell(13, 1).Range.Text = "16÷2="
$t.Cell(13, 2).Range.Text = "83÷2="
$t.Cell(13, 3).Range.Text = "70÷7="
$t.Cell(13, 4).Range.Text = "51÷2="
$t.Cell(13, 5).Range.Text = "41÷6="

$t.Cell(17, 1).Range.Text = "64÷7="
$t.Cell(17, 2).Range.Text = "56÷5="
$t.Cell(17, 3).Range.Text = "18÷2="
$t.Cell(17, 4).Range.Text = "85÷3="
$t.Cell(17, 5).Range.Text = "57÷8="
